$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 8 (shifts old rows 8,9,10,11,12 down to 9,10,11,12,13).
# Excel will automatically extend the shared formulas (B,H,I,K,L columns) that
# straddle the insertion point, and fill the new row 8 with copies of the
# formulas from the surrounding rows.
$ws.Rows.Item(8).Insert()

# New angle value for the inserted row.
$ws.Range("A8").Value = 65

# The row insert does not auto-fill the formulas into the freshly created
# row, so re-create them explicitly, matching the pattern used by every
# other data row.
$ws.Range("B8").Formula = "=A8*3.14159265/180"
$ws.Range("H8").Formula = "=`$C`$2*COS(`$B8)*`$F`$2+`$D`$2"
$ws.Range("I8").Formula = "=C`$2*SIN(`$B8)*`$F`$2-`$G`$2*(`$F`$2)^2/2+`$E`$2"
$ws.Range("K8").Formula = "=`$C`$3*COS(`$B8)*`$F`$2+`$D`$2"
$ws.Range("L8").Formula = "=F`$2*SIN(`$B8)*`$F`$2-`$G`$2*(`$F`$2)^2/2+`$E`$2"

# Update the Time (F2) input value used throughout the parabola calculations.
$ws.Range("F2").Value = 0.7

# Update the selected cell shown when the workbook is opened.
$ws.Range("L8").Select()

# Best-effort: restore on-screen window geometry recorded in the saved file
# (purely cosmetic; harmless if the host does not expose it).
$win = $wb.Windows.Item(1)
$win.Left = 1540
$win.Top = 1960
$win.Width = 19040
$win.Height = 16620

$wb.Save()
